$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text fix (4 occurrences) - typo + trailing period
$preconditionNew = "O usuário devidamente autenticado e na tela inicial de cancelar diárias."
$ws.Range("B8").Value = $preconditionNew
$ws.Range("B15").Value = $preconditionNew
$ws.Range("B23").Value = $preconditionNew
$ws.Range("B32").Value = $preconditionNew

# Fix typo "Solcitação" -> "Solicitação" in TC1 error message
$ws.Range("D10").Value = "SYSTEM Identifica que a solicitação de diária está em situação diferente de 'SOLICITADA PARA EMPENHO' ou 'SOLICITADA PARA PRESTAÇÃO DE CONTAS'.  Impede o cancelamento e exibe mensagem de erro (MSG205 - Solicitação de diária não pode ser cancelada) para o usuário."

# Add trailing period to MSG102 confirmation message (3 occurrences)
$confirmNew = "SYSTEM Exibe a mensagem (MSG102 - Confirmar cancelamento)."
$ws.Range("D17").Value = $confirmNew
$ws.Range("D25").Value = $confirmNew
$ws.Range("D34").Value = $confirmNew

# Remove stray tab character before closing parenthesis in MSG217 message
$ws.Range("D18").Value = "SYSTEM Identifica que o usuário não informou uma justificativa para o cancelamento. Não efetiva o cancelamento e exibe mensagem de erro (MSG217 - Necessário informar uma justificativa para o cancelamento de solicitações) para o usuário."
